{"js": "// fix bug in mk_frq_cnts to include entire text node\n//\n// Inserts a brand-new, empty paragraph at the very start of the document\n// body. The paragraph carries the same \"placeholder tag\" formatting used\n// throughout this generated transcription (widowControl off, nil paragraph\n// borders, clear shading, no contextual spacing, and a paragraph-mark run\n// font of 9pt gray Courier New) but holds no text - just the trailing\n// empty run (rtl=0) that this generator always appends.\n//\n// Word's high-level insertParagraph() would normally clone the pPr/rPr of\n// the paragraph it is inserted next to (copying the first paragraph's\n// black, default-font mark instead of the gray Courier New one the diff\n// calls for), so the exact target formatting is supplied directly via a\n// FlatOPC insertOoxml() at the start of the body.\nconst newParagraphXml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:pPr>' +\n      '<w:widowControl w:val=\"0\"/>' +\n      '<w:pBdr>' +\n        '<w:top w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:left w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:bottom w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:right w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:between w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n      '</w:pBdr>' +\n      '<w:shd w:fill=\"auto\" w:val=\"clear\"/>' +\n      '<w:contextualSpacing w:val=\"0\"/>' +\n      '<w:rPr>' +\n        '<w:rFonts w:ascii=\"Courier New\" w:cs=\"Courier New\" w:eastAsia=\"Courier New\" w:hAnsi=\"Courier New\"/>' +\n        '<w:color w:val=\"a9a9a9\"/>' +\n        '<w:sz w:val=\"18\"/>' +\n        '<w:szCs w:val=\"18\"/>' +\n      '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rtl w:val=\"0\"/>' +\n      '</w:rPr>' +\n    '</w:r>' +\n  '</w:p>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n      '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n          '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + newParagraphXml + '<w:sectPr/></w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nconst body = context.document.body;\nbody.insertOoxml(flatOpc, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# fix bug in mk_frq_cnts to include entire text node\n#\n# Inserts a brand-new, empty paragraph at the very start of the document\n# body. The paragraph carries the same \"placeholder tag\" formatting used\n# throughout this generated transcription (widowControl off, nil paragraph\n# borders, clear shading, no contextual spacing, and a paragraph-mark run\n# font of 9pt gray Courier New) but holds no text - just the trailing\n# empty run (rtl=0) that this generator always appends.\n#\n# Range.InsertParagraphBefore() would clone the pPr/rPr of the neighboring\n# paragraph (the first paragraph's black, default-font mark, not the gray\n# Courier New one the diff calls for), so the exact target formatting is\n# supplied directly via Range.InsertXML() on a collapsed range at the very\n# start of the document (position 0).\n\n$d = $word.ActiveDocument\n\n$newParagraphXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:pPr>' +\n      '<w:widowControl w:val=\"0\"/>' +\n      '<w:pBdr>' +\n        '<w:top w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:left w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:bottom w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:right w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n        '<w:between w:space=\"0\" w:sz=\"0\" w:val=\"nil\"/>' +\n      '</w:pBdr>' +\n      '<w:shd w:fill=\"auto\" w:val=\"clear\"/>' +\n      '<w:contextualSpacing w:val=\"0\"/>' +\n      '<w:rPr>' +\n        '<w:rFonts w:ascii=\"Courier New\" w:cs=\"Courier New\" w:eastAsia=\"Courier New\" w:hAnsi=\"Courier New\"/>' +\n        '<w:color w:val=\"a9a9a9\"/>' +\n        '<w:sz w:val=\"18\"/>' +\n        '<w:szCs w:val=\"18\"/>' +\n      '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rtl w:val=\"0\"/>' +\n      '</w:rPr>' +\n    '</w:r>' +\n  '</w:p>'\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n      '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n          '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' + $newParagraphXml + '<w:sectPr/></w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n$r = $d.Range(0, 0)\n$r.InsertXML($flatOpc)\n"}
